$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CRMAccuracyData")

$newRow = 48

$ws.Cells.Item($newRow, 1).Value = 20211018
$ws.Cells.Item($newRow, 2).Value = 2221.73708567331
$ws.Cells.Item($newRow, 3).Value = 2225.4699999999998
$ws.Cells.Item($newRow, 4).Formula = "=100*(B48-C48)/C48"
$ws.Cells.Item($newRow, 5).Value = 181
$ws.Cells.Item($newRow, 6).Value = "CRM OPENED 20210721"

$ws.Range("A49").Select()
